# Update cryptos list - GitHub Actions style refresh of price/volume figures,
# plus a couple of row-order swaps (MXToken/RenderToken and EnergySwap/PaxDollar).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay plain text
# (Price column values like "307.91" look numeric to Excel and would
# otherwise get auto-converted into a real number). We briefly mark the
# cell as Text, assign, then clear the formatting again so no stray style
# index gets left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Simple per-row updates: Price (D) and/or Volume(1h) (E) ---
# Each entry: row, newPrice (or $null to leave untouched), newVolume (or $null to leave untouched)
$updates = @(
    @{ Row = 2;  D = "27.252.72";     E = "  +0.42%  " },
    @{ Row = 3;  D = "1.908.24";      E = "  +0.32%  " },
    @{ Row = 4;  D = $null;           E = "  +0.16%  " },
    @{ Row = 5;  D = "307.91";        E = "  +0.66%  " },
    @{ Row = 6;  D = $null;           E = "  +0.21%  " },
    @{ Row = 7;  D = "0.5252";        E = $null },
    @{ Row = 8;  D = "0.3821";        E = "  +1.56%  " },
    @{ Row = 9;  D = "0.07310";       E = "  +0.87%  " },
    @{ Row = 10; D = "21.59";         E = "  +2.29%  " },
    @{ Row = 11; D = "0.9065";        E = "  +0.62%  " },
    @{ Row = 12; D = "0.08104";       E = "  -4.44%  " },
    @{ Row = 13; D = "96.21";         E = "  +1.18%  " },
    @{ Row = 14; D = "5.375";         E = "  +1.64%  " },
    @{ Row = 15; D = "1.774.37";      E = "  -6.86%  " },
    @{ Row = 16; D = "1.001";         E = "  +0.11%  " },
    @{ Row = 17; D = "0.000008683";   E = "  +0.69%  " },
    @{ Row = 18; D = "14.77";         E = "  +1.71%  " },
    @{ Row = 19; D = $null;           E = "  +0.17%  " },
    @{ Row = 20; D = "27.283.85";     E = "  +0.38%  " },
    @{ Row = 21; D = "5.122";         E = "  +1.17%  " },
    @{ Row = 22; D = "10.83";         E = "  +2.16%  " },
    @{ Row = 23; D = "6.489";         E = "  +1.02%  " },
    @{ Row = 24; D = "2.349";         E = "  +2.91%  " },
    @{ Row = 25; D = "150.23";        E = "  +2.09%  " },
    @{ Row = 26; D = "18.26";         E = "  +0.35%  " },
    @{ Row = 27; D = "1.742";         E = "  -0.57%  " },
    @{ Row = 28; D = "117.01";        E = "  +1.87%  " },
    @{ Row = 29; D = "4.858";         E = "  +1.07%  " },
    @{ Row = 30; D = "4.884";         E = "  -0.13%  " },
    @{ Row = 31; D = "0.09231";       E = "  -0.22%  " },
    @{ Row = 32; D = "0.8220";        E = "  +1.87%  " },
    @{ Row = 33; D = "0.05082";       E = $null },
    @{ Row = 34; D = "1.236";         E = "  +0.07%  " },
    @{ Row = 35; D = "2.996";         E = "  +1.84%  " },
    @{ Row = 38; D = "0.5760";        E = "  +0.95%  " },
    @{ Row = 39; D = "0.02008";       E = "  +0.67%  " },
    @{ Row = 40; D = "1.085";         E = "  +0.95%  " },
    @{ Row = 41; D = "9.071";         E = "  +0.54%  " },
    @{ Row = 42; D = "6.618";         E = "  -0.07%  " },
    @{ Row = 43; D = "117.05";        E = "  +0.66%  " },
    @{ Row = 44; D = "0.1525";        E = "  +0.68%  " },
    @{ Row = 45; D = "0.4948";        E = "  +1.94%  " },
    @{ Row = 48; D = "1.645";         E = "  +1.81%  " },
    @{ Row = 49; D = "38.70";         E = "  +3.33%  " },
    @{ Row = 50; D = $null;           E = "  +0.45%  " },
    @{ Row = 51; D = "0.05970";       E = "  +0.47%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextValue $ws.Range("D$r") $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$r").Value = $u.E
    }
}

# --- Rows 36/37: MXToken and RenderToken swap places (with refreshed figures) ---
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D36") "2.756"
$ws.Range("E36").Value = "  +4.93%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D37") "3.366"
$ws.Range("E37").Value = "  -2.17%  "

# --- Rows 46/47: EnergySwap and PaxDollar swap places (with refreshed figures) ---
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D46") "1.002"
$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "10.15"
$ws.Range("E47").Value = "  +0.38%  "
